# Auto update Excel log
# Adds two new sheets ("Proximity" and "Camera") with sensor log rows,
# placed after the existing "ALERTS" and "mmWave" sheets.

$wb = $excel.ActiveWorkbook

$headers = @("Date", "Timestamp", "Hour", "Location", "Value", "Status")

# ---- Proximity sheet ------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$proximity = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$proximity.Name = "Proximity"

$proximityData = @(
    @("2026-01-30", "14:47:06", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-30", "14:47:25", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-01-30", "14:47:27", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

# Ensure all cells used (including date-looking text) are kept as plain text,
# matching the source log format instead of being auto-converted to dates.
$proximity.Range("A1:F" + ($proximityData.Length + 1)).NumberFormat = "@"

for ($col = 0; $col -lt $headers.Length; $col++) {
    $proximity.Cells.Item(1, $col + 1).Value = $headers[$col]
}

for ($row = 0; $row -lt $proximityData.Length; $row++) {
    $rowData = $proximityData[$row]
    for ($col = 0; $col -lt $rowData.Length; $col++) {
        $proximity.Cells.Item($row + 2, $col + 1).Value = $rowData[$col]
    }
}

# ---- Camera sheet -----------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$camera = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$camera.Name = "Camera"

$cameraData = @(
    @("2026-01-30", "14:47:07", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-30", "14:47:14", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-30", "14:47:23", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-30", "14:47:28", "14:00", "Living Room Main Door", "Image Captured", "Active")
)

$camera.Range("A1:F" + ($cameraData.Length + 1)).NumberFormat = "@"

for ($col = 0; $col -lt $headers.Length; $col++) {
    $camera.Cells.Item(1, $col + 1).Value = $headers[$col]
}

for ($row = 0; $row -lt $cameraData.Length; $row++) {
    $rowData = $cameraData[$row]
    for ($col = 0; $col -lt $rowData.Length; $col++) {
        $camera.Cells.Item($row + 2, $col + 1).Value = $rowData[$col]
    }
}

Write-Host "Sheets now:" $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name
}
